# 陳明文 2011-11-23 財產申報表 — split "其他有價證券" into
#   具有相當價值之財產 (antiques / other valuable property)
#   保險 (insurance)               <- brand new sheet
#   債務 (debts)                   <- re-numbered (sheetId 6 -> 7)
#
# Strategy:
#  - Delete the old "債務" sheet first so its sheetId (6) is freed up;
#    this lets the newly inserted "保險" sheet claim sheetId 6 and the
#    recreated "債務" sheet claim sheetId 7, matching the target layout.
#  - Rename sheet5 "其他有價證券" -> "具有相當價值之財產" and rewrite its data.
#  - Insert the new "保險" sheet right after sheet5 and fill it in.
#  - Insert a fresh "債務" sheet right after "保險" and re-populate it with
#    the same debt data the original sheet held (only the shared-string
#    slots shift; the values themselves are unchanged).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 0: remember header/data style donors (plain cells with the
# existing bold-header style "1" / plain-data style "2") before we
# start reshaping things, and drop + rebuild the debt sheet so the
# sheetId sequence comes out 5 / 6 / 7.
# ---------------------------------------------------------------------
$wsDebtOld = $wb.Worksheets.Item("債務")
$wsDebtOld.Delete()

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "具有相當價值之財產"

$wsIns = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws5)
$wsIns.Name = "保險"

$wsDebt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsIns)
$wsDebt.Name = "債務"

# ---------------------------------------------------------------------
# Step 1: 具有相當價值之財產 (sheet5) — collapse the old 11-row table into
# a 2-row (header + 1 data row) table spanning columns A:L.
# ---------------------------------------------------------------------

# Extend the existing header / data-row formatting (style "1" / "2")
# out to the new columns H:L before touching any values.
$ws5.Range("D1").Copy()
$ws5.Range("H1:L1").PasteSpecial(-4122)
$ws5.Range("D2").Copy()
$ws5.Range("H2:L2").PasteSpecial(-4122)

# Drop the now-obsolete rows 3-11 entirely.
$ws5.Rows("3:11").Clear()

$ws5.Range("B1").Value = "name"
$ws5.Range("C1").Value = "quantity"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "property_category"
$ws5.Range("G1").Value = "category"
$ws5.Range("H1").Value = "date"
$ws5.Range("I1").Value = "legislator_name"
$ws5.Range("J1").Value = "legislator_id"
$ws5.Range("K1").Value = "source_file"
$ws5.Range("L1").Value = "index"

$ws5.Range("A2").Value = 120
$ws5.Range("B2").Value = "朱銘雕刻"
$ws5.Range("C2").Value = 1
$ws5.Range("D2").Value = "陳明文"
$ws5.Range("E2").Value = 500000
$ws5.Range("F2").Value = "otherbonds"
$ws5.Range("G2").Value = "normal"
$ws5.Range("H2").NumberFormat = "@"
$ws5.Range("H2").Value = "2011-11-23"
$ws5.Range("I2").Value = "陳明文"
$ws5.Range("J2").Value = 828
$ws5.Range("K2").Value = "tmpf421"
$ws5.Range("L2").Value = 120

# ---------------------------------------------------------------------
# Step 2: 保險 (sheet6, brand new) — header + 6 insurance policy rows,
# columns A:E.
# ---------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("土地")
$ws1.Range("B1").Copy()
$wsIns.Range("B1:E1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$wsIns.Range("A2:A7").PasteSpecial(-4122)
$ws1.Range("B2").Copy()
$wsIns.Range("B2:E7").PasteSpecial(-4122)

$wsIns.Range("B1").Value = "國泰人壽"
$wsIns.Range("C1").Value = "得意還本終身"
$wsIns.Range("D1").Value = "陳明文"
$wsIns.Range("E1").Value = "88.05.07108.05.7"

$wsIns.Range("A2").Value = 123
$wsIns.Range("B2").Value = "國泰人壽"
$wsIns.Range("C2").Value = "得意還本終身"
$wsIns.Range("D2").Value = "陳明文"
$wsIns.Range("E2").Value = "88.05.07108.05.7"

$wsIns.Range("A3").Value = 124
$wsIns.Range("B3").Value = "國泰人壽"
$wsIns.Range("C3").Value = "新富貴保本投資鏈結型保險第7期"
$wsIns.Range("D3").Value = "廖素惠"
$wsIns.Range("E3").Value = "100.06.13106.06.12"

$wsIns.Range("A4").Value = 125
$wsIns.Range("B4").Value = "國泰人壽"
$wsIns.Range("C4").Value = "創世變額萬能壽險"
$wsIns.Range("D4").Value = "廖素惠"
$wsIns.Range("E4").Value = "94.12.30100.12.30"

$wsIns.Range("A5").Value = 126
$wsIns.Range("B5").Value = "國泰人壽"
$wsIns.Range("C5").Value = "創世變額萬能壽險"
$wsIns.Range("D5").Value = "廖素惠"
$wsIns.Range("E5").Value = "94.12.30100.12.30(被保人:陳〇廷）"

$wsIns.Range("A6").Value = 127
$wsIns.Range("B6").Value = "富邦人壽"
$wsIns.Range("C6").Value = "安泰增額養老壽險"
$wsIns.Range("D6").Value = "陳明文"
$wsIns.Range("E6").Value = "86.04.10101.04.10"

$wsIns.Range("A7").Value = 128
$wsIns.Range("B7").Value = "台灣人壽"
$wsIns.Range("C7").Value = "台灣人壽富利長紅利率變動型年金保險"
$wsIns.Range("D7").Value = "陳〇廷"
$wsIns.Range("E7").Value = "97.05.19"

# ---------------------------------------------------------------------
# Step 3: 債務 (sheet7, recreated) — same debt entry the original sheet
# held, columns A:G.
# ---------------------------------------------------------------------

$ws1.Range("B1").Copy()
$wsDebt.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$wsDebt.Range("A2").PasteSpecial(-4122)
$ws1.Range("B2").Copy()
$wsDebt.Range("B2:G2").PasteSpecial(-4122)

$wsDebt.Range("B1").Value = "一般借款"
$wsDebt.Range("C1").Value = "陳明文"
$wsDebt.Range("D1").Value = "嘉義縣朴子市農會本會嘉義縣朴子市山通路"
$wsDebt.Range("E1").Value = "7425233"
$wsDebt.Range("F1").Value = "98年12月20日"
$wsDebt.Range("G1").Value = "償還擔保債務&房屋整修&新建房屋"

$wsDebt.Range("A2").Value = 138
$wsDebt.Range("B2").Value = "一般借款"
$wsDebt.Range("C2").Value = "陳明文"
$wsDebt.Range("D2").Value = "嘉義縣朴子市農會本會嘉義縣朴子市山通路"
$wsDebt.Range("E2").Value = "7425233"
$wsDebt.Range("F2").Value = "98年12月20日"
$wsDebt.Range("G2").Value = "償還擔保債務&房屋整修&新建房屋"
